# Added User Site Setting wheel navigation section
# (new block of rows at the bottom of the Objects_Navigation sheet,
#  mirroring the layout/style of the other "section header + rows" blocks
#  already on the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Navigation")

# --- Section header row 30: "User Site Setting wheel navigation" ---
# Copy the formatting (fill/font/merge look) from an existing section header
# row (A2:G2) onto the new row, then set the header text and merge A30:G30
# exactly like the other section headers on this sheet.
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A30:G30").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A30").Value = "User Site Setting wheel navigation"
$ws.Range("A30:G30").Merge() | Out-Null

# --- Row 31: Settings Wheel ---
$ws.Range("B31").Value = "Settings Wheel"
$ws.Range("C31").Value = "xpath"
$ws.Range("D31").Value = "button"
$ws.Range("E31").Value = ".//*[@id='siteactiontd']"

# --- Row 32: Site contents ---
$ws.Range("B32").Value = "Site contents"
$ws.Range("C32").Value = "xpath"
$ws.Range("D32").Value = "link"
$ws.Range("E32").Value = "//*[ text()='Site contents']"

# --- Row 33: add an app ---
$ws.Range("B33").Value = "add an app"
$ws.Range("C33").Value = "xpath"
$ws.Range("D33").Value = "button"
$ws.Range("E33").Value = ".//*[@id='apptile-appadd']/div[1]/a"

# New rows use the same "id,name,xpath" dropdown validation on the
# LocatorType column as the rest of the sheet.
$ws.Range("C31:C33").Validation.Add(3, 1, 1, '"id,name,xpath"') | Out-Null

# Leave the selection where the author left it when they saved.
$ws.Activate() | Out-Null
$ws.Range("J30").Select() | Out-Null
